# Update cryptos list figures (Price / Volume(1h)) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.514.12"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").Value = "'1.831.97"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'313.93"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").Value = "'1.001"

$ws.Range("D7").Value = "'0.4298"
$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("D8").Value = "'0.3661"
$ws.Range("E8").Value = "  +0.33%  "

$ws.Range("D9").Value = "'0.07272"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("D10").Value = "'0.8691"
$ws.Range("E10").Value = "  -0.77%  "

$ws.Range("D11").Value = "'20.67"
$ws.Range("E11").Value = "  -0.43%  "

$ws.Range("D12").Value = "'1.816.33"
$ws.Range("E12").Value = "  +1.17%  "

$ws.Range("D13").Value = "'5.415"
$ws.Range("E13").Value = "  +1.39%  "

$ws.Range("D14").Value = "'6.541"
$ws.Range("E14").Value = "  +0.33%  "

$ws.Range("D15").Value = "'0.06929"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").Value = "'80.64"
$ws.Range("E17").Value = "  +0.84%  "

$ws.Range("D18").Value = "'0.000008905"
$ws.Range("E18").Value = "  -0.87%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.22%  "

$ws.Range("D20").Value = "'15.43"
$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").Value = "'27.750.99"
$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("D22").Value = "'5.149"
$ws.Range("E22").Value = "  +3.41%  "

$ws.Range("D23").Value = "'10.83"
$ws.Range("E23").Value = "  +4.17%  "

$ws.Range("D24").Value = "'2.123.22"
$ws.Range("E24").Value = "  +4.49%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").Value = "'154.54"
$ws.Range("E26").Value = "  -0.96%  "

$ws.Range("D27").Value = "'18.84"
$ws.Range("E27").Value = "  +1.17%  "

$ws.Range("D28").Value = "'5.157"
$ws.Range("E28").Value = "  -1.76%  "

$ws.Range("D29").Value = "'114.15"
$ws.Range("E29").Value = "  -5.07%  "

$ws.Range("D30").Value = "'1.833"
$ws.Range("E30").Value = "  -1.30%  "

$ws.Range("D31").Value = "'0.08891"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").Value = "'0.7561"
$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").Value = "'2.985"
$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("D34").Value = "'4.548"
$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("D35").Value = "'1.138"
$ws.Range("E35").Value = "  +1.38%  "

$ws.Range("D36").Value = "'1.001"

$ws.Range("D37").Value = "'1.092"
$ws.Range("E37").Value = "  -1.23%  "

$ws.Range("D38").Value = "'0.05316"
$ws.Range("E38").Value = "  -1.76%  "

$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("D40").Value = "'2.797"
$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("D41").Value = "'0.1669"
$ws.Range("E41").Value = "  +0.93%  "

$ws.Range("D42").Value = "'0.5071"
$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").Value = "'6.619"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").Value = "'8.398"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").Value = "'10.42"
$ws.Range("E45").Value = "  +1.53%  "

$ws.Range("D46").Value = "'106.15"
$ws.Range("E46").Value = "  +1.85%  "

$ws.Range("E47").Value = "  -0.58%  "

$ws.Range("D48").Value = "'0.4686"
$ws.Range("E48").Value = "  +0.45%  "

$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("D50").Value = "'1.609"
$ws.Range("E50").Value = "  -0.85%  "

$ws.Range("D51").Value = "'64.19"
$ws.Range("E51").Value = "  -0.04%  "
